$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range('D2')
$r.NumberFormat = '@'
$r.Value = '62.927.59'
$r.Style = 'Normal'
$ws.Range('E2').Value = '  +2.03%  '
$ws.Range('E3').Value = '  +1.16%  '
$ws.Range('E4').Value = '  +0.02%  '
$r = $ws.Range('D5')
$r.NumberFormat = '@'
$r.Value = '593.51'
$r.Style = 'Normal'
$ws.Range('E5').Value = '  -0.45%  '
$r = $ws.Range('D6')
$r.NumberFormat = '@'
$r.Value = '154.55'
$r.Style = 'Normal'
$ws.Range('E6').Value = '  +6.93%  '
$r = $ws.Range('D7')
$r.NumberFormat = '@'
$r.Value = '0.999'
$r.Style = 'Normal'
$ws.Range('E7').Value = '  -0.15%  '
$r = $ws.Range('D8')
$r.NumberFormat = '@'
$r.Value = '3.027.14'
$r.Style = 'Normal'
$ws.Range('E8').Value = '  +1.01%  '
$r = $ws.Range('D10')
$r.NumberFormat = '@'
$r.Value = '6.45'
$r.Style = 'Normal'
$ws.Range('E10').Value = '  +8.94%  '
$ws.Range('E11').Value = '  +1.89%  '
$r = $ws.Range('D12')
$r.NumberFormat = '@'
$r.Value = '0.467'
$r.Style = 'Normal'
$ws.Range('E12').Value = '  +1.21%  '
$r = $ws.Range('D13')
$r.NumberFormat = '@'
$r.Value = '0.0000235'
$r.Style = 'Normal'
$ws.Range('E13').Value = '  +2.14%  '
$r = $ws.Range('D14')
$r.NumberFormat = '@'
$r.Value = '35.62'
$r.Style = 'Normal'
$ws.Range('E14').Value = '  +3.78%  '
$ws.Range('E15').Value = '  +2.27%  '
$r = $ws.Range('D16')
$r.NumberFormat = '@'
$r.Value = '3.536.08'
$r.Style = 'Normal'
$ws.Range('E16').Value = '  +1.15%  '
$ws.Range('E17').Value = '  +0.85%  '
$r = $ws.Range('D18')
$r.NumberFormat = '@'
$r.Value = '62.897.36'
$r.Style = 'Normal'
$ws.Range('E18').Value = '  +2.03%  '
$r = $ws.Range('D19')
$r.NumberFormat = '@'
$r.Value = '3.032.01'
$r.Style = 'Normal'
$ws.Range('E19').Value = '  +0.98%  '
$r = $ws.Range('D20')
$r.NumberFormat = '@'
$r.Value = '453.16'
$r.Style = 'Normal'
$ws.Range('E20').Value = '  -0.29%  '
$r = $ws.Range('D21')
$r.NumberFormat = '@'
$r.Value = '14.31'
$r.Style = 'Normal'
$ws.Range('E21').Value = '  +2.15%  '
$r = $ws.Range('D22')
$r.NumberFormat = '@'
$r.Value = '0.698'
$r.Style = 'Normal'
$ws.Range('E22').Value = '  +1.37%  '
$r = $ws.Range('D23')
$r.NumberFormat = '@'
$r.Value = '7.49'
$r.Style = 'Normal'
$ws.Range('E23').Value = '  +1.80%  '
$r = $ws.Range('D24')
$r.NumberFormat = '@'
$r.Value = '83.12'
$r.Style = 'Normal'
$ws.Range('E24').Value = '  +1.27%  '
$ws.Range('B25').Value = 'RenderToken'
$ws.Range('C25').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$r = $ws.Range('D25')
$r.NumberFormat = '@'
$r.Value = '11.21'
$r.Style = 'Normal'
$ws.Range('E25').Value = '  +6.40%  '
$ws.Range('B26').Value = 'Fetch.AI'
$ws.Range('C26').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$r = $ws.Range('D26')
$r.NumberFormat = '@'
$r.Value = '2.32'
$r.Style = 'Normal'
$ws.Range('E26').Value = '  +3.36%  '
$r = $ws.Range('D27')
$r.NumberFormat = '@'
$r.Value = '12.38'
$r.Style = 'Normal'
$ws.Range('E27').Value = '  +2.65%  '
$r = $ws.Range('D29')
$r.NumberFormat = '@'
$r.Value = '7.46'
$r.Style = 'Normal'
$ws.Range('E29').Value = '  +3.96%  '
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('E31').Value = '  -0.04%  '
$ws.Range('E32').Value = '  +6.03%  '
$ws.Range('E33').Value = '  +0.25%  '
$r = $ws.Range('D34')
$r.NumberFormat = '@'
$r.Value = '0.111'
$r.Style = 'Normal'
$ws.Range('E34').Value = '  +2.04%  '
$r = $ws.Range('D35')
$r.NumberFormat = '@'
$r.Value = '0.0₃0872'
$r.Style = 'Normal'
$ws.Range('E35').Value = '  +4.60%  '
$ws.Range('E36').Value = '  +2.00%  '
$r = $ws.Range('D37')
$r.NumberFormat = '@'
$r.Value = '5.94'
$r.Style = 'Normal'
$ws.Range('E37').Value = '  +2.88%  '
$ws.Range('E38').Value = '  +10.21%  '
$r = $ws.Range('D39')
$r.NumberFormat = '@'
$r.Value = '2.12'
$r.Style = 'Normal'
$ws.Range('E39').Value = '  +2.42%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$r = $ws.Range('D40')
$r.NumberFormat = '@'
$r.Value = '50.66'
$r.Style = 'Normal'
$ws.Range('E40').Value = '  +0.76%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$r = $ws.Range('D41')
$r.NumberFormat = '@'
$r.Value = '0.128'
$r.Style = 'Normal'
$ws.Range('E41').Value = '  +4.50%  '
$r = $ws.Range('D42')
$r.NumberFormat = '@'
$r.Value = '9.08'
$r.Style = 'Normal'
$ws.Range('E42').Value = '  -2.05%  '
$r = $ws.Range('D43')
$r.NumberFormat = '@'
$r.Value = '0.308'
$r.Style = 'Normal'
$ws.Range('E43').Value = '  +14.20%  '
$r = $ws.Range('D44')
$r.NumberFormat = '@'
$r.Value = '41.83'
$r.Style = 'Normal'
$ws.Range('E44').Value = '  +6.05%  '
$r = $ws.Range('D45')
$r.NumberFormat = '@'
$r.Value = '396.50'
$r.Style = 'Normal'
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('E46').Value = '  +1.23%  '
$r = $ws.Range('D47')
$r.NumberFormat = '@'
$r.Value = '2.731.43'
$r.Style = 'Normal'
$ws.Range('E47').Value = '  +0.29%  '
$r = $ws.Range('D48')
$r.NumberFormat = '@'
$r.Value = '132.46'
$r.Style = 'Normal'
$ws.Range('E48').Value = '  -0.62%  '
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('E50').Value = '  +3.54%  '
$r = $ws.Range('D51')
$r.NumberFormat = '@'
$r.Value = '24.47'
$r.Style = 'Normal'
$ws.Range('E51').Value = '  +3.92%  '
